$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Names for the new first column (row 2..12), keyed by row number
$names = @{
    2  = "Hayle Estuary"
    3  = "Gannel Estuary"
    4  = "Camel Estuary"
    5  = "Bridgwater Bay"
    6  = "Severn Estuary"
    7  = "Axe Estuary"
    8  = "Otter Estuary"
    9  = "Exe Estuary"
    10 = "Teign Estuary"
    11 = "Dart Estuary"
    12 = "Salcombe & Kingsbridge Estuary"
}

# Shift the existing three columns (A,B,C) one column to the right (-> B,C,D),
# reading all source values first so the write order cannot clobber data we
# still need.
for ($r = 1; $r -le 12; $r++) {
    $valA = $ws.Cells.Item($r, 1).Value()
    $valB = $ws.Cells.Item($r, 2).Value()
    $valC = $ws.Cells.Item($r, 3).Value()

    $ws.Cells.Item($r, 4).Value = $valC
    $ws.Cells.Item($r, 3).Value = $valB
    $ws.Cells.Item($r, 2).Value = $valA
}

# New column D (old column C) holds numeric values for rows 2-12; the moved
# cells are brand new so they need the 2-decimal number format re-applied.
$ws.Range("D2:D12").NumberFormat = "0.00"

# Populate the new first column with the row names.
$ws.Cells.Item(1, 1).Value = "Name"
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = $names[$r]
}

# New column A width.
$ws.Columns("A:A").ColumnWidth = 14.5

# Update the selected cell shown in the sheet view.
$ws.Range("N20").Select()
